$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.735
$ws.Range("D3").Value = -7.827
$ws.Range("E19").Value = 16.441
$ws.Range("A21").Value = -20.149
$ws.Range("A23").Value = -20.142
$ws.Range("D24").Value = -7.527000000000001
$ws.Range("E24").Value = 17
$ws.Range("A25").Value = -21.78900000000001
$ws.Range("C27").Value = -13.147
$ws.Range("E30").Value = 16.394
$ws.Range("C31").Value = -12.85
$ws.Range("E31").Value = 16.319
$ws.Range("E33").Value = 17.253
$ws.Range("C39").Value = -12.847
$ws.Range("C48").Value = -11.181
$ws.Range("C51").Value = -11.506
$ws.Range("C52").Value = -11.496
$ws.Range("A53").Value = -21.794
$ws.Range("C55").Value = -13.644
$ws.Range("E55").Value = 16.409
$ws.Range("C56").Value = -12.996
$ws.Range("A57").Value = -21.942
$ws.Range("C57").Value = -12.624
$ws.Range("D57").Value = -8.269000000000002
$ws.Range("A59").Value = -22.358
$ws.Range("D61").Value = -7.708
$ws.Range("E65").Value = 17.36
$ws.Range("A69").Value = -21.703
$ws.Range("D70").Value = -7.419999999999999
$ws.Range("E70").Value = 17.612
$ws.Range("C73").Value = -12.668
$ws.Range("E75").Value = 16.76
$ws.Range("A79").Value = -20.901
$ws.Range("A83").Value = -21.938
$ws.Range("E83").Value = 16.72
$ws.Range("D86").Value = -8.186
$ws.Range("C89").Value = -11.202
$ws.Range("C90").Value = -12.942
$ws.Range("A93").Value = -21.481
$ws.Range("E96").Value = 16.454
$ws.Range("E97").Value = 16.948
$ws.Range("D98").Value = -8.397
$ws.Range("D100").Value = -8.361999999999998
$ws.Range("D102").Value = -7.794000000000001
